$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1066.575
$ws.Range("I40").Value = 1070.7949
$ws.Range("J40").Value = 902
$ws.Range("K40").Value = 1070.7949
$ws.Range("L40").Value = 902
$ws.Range("M40").Value = -895.7949000000001
$ws.Range("N40").Value = -1252

$ws.Range("H62").Value = 1800
$ws.Range("I62").Value = 1777.7778
$ws.Range("K62").Value = 1777.7778
$ws.Range("M62").Value = -1153.7778

$ws.Range("H65").Value = 1800
$ws.Range("I65").Value = 1777.7778
$ws.Range("K65").Value = 8888.889000000001
$ws.Range("M65").Value = -5768.889000000001

$ws.Range("H76").Value = 3023.75
$ws.Range("I76").Value = 2500
$ws.Range("J76").Value = 3198.3333
$ws.Range("K76").Value = 2500
$ws.Range("L76").Value = 3198.3333
$ws.Range("M76").Value = -2185
$ws.Range("N76").Value = -3828.3333

$ws.Range("H79").Value = 3023.75
$ws.Range("I79").Value = 2500
$ws.Range("J79").Value = 3198.3333
$ws.Range("K79").Value = 2500
$ws.Range("L79").Value = 3198.3333
$ws.Range("M79").Value = -1408
$ws.Range("N79").Value = -5382.3333

$ws.Range("H88").Value = 2928.2856
$ws.Range("I88").Value = 4333
$ws.Range("K88").Value = 4333
$ws.Range("M88").Value = -3927

$ws.Range("H91").Value = 2928.2856
$ws.Range("I91").Value = 4333
$ws.Range("K91").Value = 4333
$ws.Range("M91").Value = -2929

$ws.Range("H107").Value = 611.7778
$ws.Range("I107").Value = 344.5
$ws.Range("J107").Value = 2750
$ws.Range("K107").Value = 344.5
$ws.Range("L107").Value = 2750
$ws.Range("M107").Value = 1575.5
$ws.Range("N107").Value = -6590

$ws.Range("H132").Value = 1214.7091
$ws.Range("I132").Value = 996.93616
$ws.Range("J132").Value = 2494.125
$ws.Range("K132").Value = 2990.80848
$ws.Range("L132").Value = 7482.375
$ws.Range("M132").Value = -460.8084799999997
$ws.Range("N132").Value = -12542.375

$ws.Range("H137").Value = 22433.766
$ws.Range("I137").Value = 975.36365
$ws.Range("J137").Value = 73014.28999999999
$ws.Range("K137").Value = 2926.09095
$ws.Range("L137").Value = 219042.87
$ws.Range("M137").Value = -376.0909499999998
$ws.Range("N137").Value = -224142.87

$ws.Range("H138").Value = 2045.7284
$ws.Range("I138").Value = 1819.2941
$ws.Range("J138").Value = 2430.6667
$ws.Range("K138").Value = 5457.8823
$ws.Range("L138").Value = 7292.000100000001
$ws.Range("M138").Value = -317.8823000000002
$ws.Range("N138").Value = -17572.0001

$ws.Range("H141").Value = 2040.5186
$ws.Range("I141").Value = 1025.52
$ws.Range("J141").Value = 14728
$ws.Range("K141").Value = 3076.56
$ws.Range("L141").Value = 44184
$ws.Range("M141").Value = 2103.44
$ws.Range("N141").Value = -54544

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4372.127
$ws.Range("I32").Value = 3941.7869
$ws.Range("J32").Value = 17497.5
$ws.Range("K32").Value = 3941.7869
$ws.Range("L32").Value = 17497.5
$ws.Range("M32").Value = -3654.7869
$ws.Range("N32").Value = -18071.5

$ws.Range("H74").Value = 1966.6538
$ws.Range("I74").Value = 1375.7858
$ws.Range("J74").Value = 2656
$ws.Range("K74").Value = 1375.7858
$ws.Range("L74").Value = 2656
$ws.Range("M74").Value = -501.7858000000001
$ws.Range("N74").Value = -4404

$ws.Range("H77").Value = 1966.6538
$ws.Range("I77").Value = 1375.7858
$ws.Range("J77").Value = 2656
$ws.Range("K77").Value = 6878.929
$ws.Range("L77").Value = 13280
$ws.Range("M77").Value = -2510.929
$ws.Range("N77").Value = -22016

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5337.6895
$ws.Range("I134").Value = 5559.72
$ws.Range("K134").Value = 16679.16
$ws.Range("M134").Value = -14144.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1008.2
$ws.Range("I16").Value = 922.4
$ws.Range("J16").Value = 1179.8
$ws.Range("K16").Value = 922.4
$ws.Range("L16").Value = 1179.8
$ws.Range("M16").Value = -635.4
$ws.Range("N16").Value = -1753.8

$ws.Range("H31").Value = 2080.394
$ws.Range("I31").Value = 1467
$ws.Range("J31").Value = 3307.182
$ws.Range("K31").Value = 1467
$ws.Range("L31").Value = 3307.182
$ws.Range("M31").Value = -1172
$ws.Range("N31").Value = -3897.182

$ws.Range("H34").Value = 2080.394
$ws.Range("I34").Value = 1467
$ws.Range("J34").Value = 3307.182
$ws.Range("K34").Value = 1467
$ws.Range("L34").Value = 3307.182
$ws.Range("M34").Value = -1265
$ws.Range("N34").Value = -3711.182

$ws.Range("H105").Value = 1270.3
$ws.Range("I105").Value = 1270.3
$ws.Range("K105").Value = 1270.3
$ws.Range("M105").Value = 476.7

$ws.Range("H113").Value = 1008.2
$ws.Range("I113").Value = 922.4
$ws.Range("J113").Value = 1179.8
$ws.Range("K113").Value = 922.4
$ws.Range("L113").Value = 1179.8
$ws.Range("M113").Value = 1247.6
$ws.Range("N113").Value = -5519.8

$ws.Range("H132").Value = 1915.625
$ws.Range("I132").Value = 1257.5807
$ws.Range("J132").Value = 4182.222
$ws.Range("K132").Value = 3772.7421
$ws.Range("L132").Value = 12546.666
$ws.Range("M132").Value = -1242.7421
$ws.Range("N132").Value = -17606.666

$ws.Range("H141").Value = 81727
$ws.Range("J141").Value = 84969.336
$ws.Range("L141").Value = 84969.336
$ws.Range("N141").Value = -95329.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6573.8423
$ws.Range("I56").Value = 6573.8423
$ws.Range("K56").Value = 6573.8423
$ws.Range("M56").Value = -6043.8423

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1856.4
$ws.Range("I102").Value = 1684.8889
$ws.Range("K102").Value = 1684.8889
$ws.Range("M102").Value = -62.88889999999992

$ws.Range("H113").Value = 1279.4546
$ws.Range("I113").Value = 1213.4286
$ws.Range("K113").Value = 1213.4286
$ws.Range("M113").Value = 956.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 41999
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 10000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -10460

$ws.Range("H40").Value = 6169.6333
$ws.Range("I40").Value = 6832.048
$ws.Range("J40").Value = 4624
$ws.Range("K40").Value = 6832.048
$ws.Range("L40").Value = 4624
$ws.Range("M40").Value = -6696.048
$ws.Range("N40").Value = -4896

$ws.Range("H61").Value = 3612.25
$ws.Range("I61").Value = 3474.5
$ws.Range("J61").Value = 3750
$ws.Range("K61").Value = 3474.5
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -3272.5
$ws.Range("N61").Value = -4154

$ws.Range("H68").Value = 1912.3077
$ws.Range("I68").Value = 1459.3334
$ws.Range("J68").Value = 2300.5715
$ws.Range("K68").Value = 1459.3334
$ws.Range("L68").Value = 2300.5715
$ws.Range("M68").Value = -710.3334
$ws.Range("N68").Value = -3798.5715

$ws.Range("H71").Value = 1912.3077
$ws.Range("I71").Value = 1459.3334
$ws.Range("J71").Value = 2300.5715
$ws.Range("K71").Value = 7296.666999999999
$ws.Range("L71").Value = 11502.8575
$ws.Range("M71").Value = -3552.666999999999
$ws.Range("N71").Value = -18990.8575

$ws.Range("H82").Value = 1391.8823
$ws.Range("I82").Value = 1008.0909
$ws.Range("J82").Value = 2095.5
$ws.Range("K82").Value = 1008.0909
$ws.Range("L82").Value = 2095.5
$ws.Range("M82").Value = -647.0909
$ws.Range("N82").Value = -2817.5

$ws.Range("H85").Value = 1391.8823
$ws.Range("I85").Value = 1008.0909
$ws.Range("J85").Value = 2095.5
$ws.Range("K85").Value = 1008.0909
$ws.Range("L85").Value = 2095.5
$ws.Range("M85").Value = 239.9091
$ws.Range("N85").Value = -4591.5

$ws.Range("H113").Value = 3612.25
$ws.Range("I113").Value = 3474.5
$ws.Range("J113").Value = 3750
$ws.Range("K113").Value = 3474.5
$ws.Range("L113").Value = 3750
$ws.Range("M113").Value = -1304.5
$ws.Range("N113").Value = -8090

$ws.Range("H132").Value = 2350.3914
$ws.Range("I132").Value = 1665.4286
$ws.Range("K132").Value = 4996.2858
$ws.Range("M132").Value = -2466.2858

$ws.Range("H136").Value = 2128.1968
$ws.Range("I136").Value = 1542.2264
$ws.Range("J136").Value = 6010.25
$ws.Range("K136").Value = 4626.6792
$ws.Range("L136").Value = 18030.75
$ws.Range("M136").Value = -2076.6792
$ws.Range("N136").Value = -23130.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 13833.333

$ws.Range("H62").Value = 7934
$ws.Range("I62").Value = 7940.8
$ws.Range("K62").Value = 7940.8
$ws.Range("M62").Value = -7316.8

$ws.Range("H65").Value = 7934
$ws.Range("I65").Value = 7940.8
$ws.Range("K65").Value = 39704
$ws.Range("M65").Value = -36584

$ws.Range("H132").Value = 1128.48
$ws.Range("I132").Value = 776.6667
$ws.Range("J132").Value = 2033.1428
$ws.Range("K132").Value = 2330.0001
$ws.Range("L132").Value = 6099.428400000001
$ws.Range("M132").Value = 199.9998999999998
$ws.Range("N132").Value = -11159.4284
